$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric ("244.52", "13", ...) must keep their original
# Text cell-type (the source sheet stores every data cell as inline text), so force the
# Text number format before writing the value - otherwise Excel auto-converts the string
# to a numeric value.
$textForceCells = @("D2", "G2", "D3", "G3", "D4", "G4", "D5", "G5", "D6", "G6", "D7", "G7", "D8", "G8", "D9", "G9", "D10", "G10", "D11", "G11", "D12", "G12", "D13", "G13", "D14", "G14", "D15", "G15", "D16", "G16", "D17", "G17", "D18", "G18", "D19", "G19", "D20", "G20", "G21", "D22", "G22", "D23", "G23", "D24", "G24", "G25", "G26", "G27", "G28", "G29", "G30", "G31", "G32", "G33", "G34", "G35", "G36", "G37", "G38", "G39", "D40", "G40", "D41", "G41", "D42", "G42", "D43", "G43", "D44", "G44", "D45", "G45", "G46", "D47", "G47", "D48", "G48", "G49", "G50", "G51")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Price column (D) + Hour column (G) updates ---
$ws.Range("D2").Value = "244.52"
$ws.Range("G2").Value = "13"
$ws.Range("D3").Value = "21.85"
$ws.Range("G3").Value = "13"
$ws.Range("D4").Value = "5.392"
$ws.Range("G4").Value = "13"
$ws.Range("D5").Value = "0.05996"
$ws.Range("G5").Value = "13"
$ws.Range("D6").Value = "3.390"
$ws.Range("G6").Value = "13"
$ws.Range("D7").Value = "0.8137"
$ws.Range("G7").Value = "13"
$ws.Range("D8").Value = "0.9531"
$ws.Range("G8").Value = "13"
$ws.Range("D9").Value = "0.1429"
$ws.Range("G9").Value = "13"
$ws.Range("D10").Value = "0.07399"
$ws.Range("G10").Value = "13"
$ws.Range("D11").Value = "0.03357"
$ws.Range("G11").Value = "13"
$ws.Range("D12").Value = "0.03060"
$ws.Range("G12").Value = "13"
$ws.Range("D13").Value = "0.09407"
$ws.Range("G13").Value = "13"
$ws.Range("D14").Value = "4.004"
$ws.Range("G14").Value = "13"
$ws.Range("D15").Value = "0.001589"
$ws.Range("G15").Value = "13"
$ws.Range("D16").Value = "0.04799"
$ws.Range("G16").Value = "13"
$ws.Range("D17").Value = "0.0005901"
$ws.Range("G17").Value = "13"
$ws.Range("D18").Value = "0.006201"
$ws.Range("G18").Value = "13"
$ws.Range("D19").Value = "0.004998"
$ws.Range("G19").Value = "13"
$ws.Range("D20").Value = "0.0009880"
$ws.Range("G20").Value = "13"
$ws.Range("G21").Value = "13"
$ws.Range("D22").Value = "3.677"
$ws.Range("G22").Value = "13"
$ws.Range("D23").Value = "6.412"
$ws.Range("G23").Value = "13"
$ws.Range("D24").Value = "2.189"
$ws.Range("G24").Value = "13"
$ws.Range("G25").Value = "13"
$ws.Range("G26").Value = "13"
$ws.Range("G27").Value = "13"
$ws.Range("G28").Value = "13"
$ws.Range("G29").Value = "13"
$ws.Range("G30").Value = "13"
$ws.Range("G31").Value = "13"
$ws.Range("G32").Value = "13"
$ws.Range("G33").Value = "13"
$ws.Range("G34").Value = "13"
$ws.Range("G35").Value = "13"
$ws.Range("G36").Value = "13"
$ws.Range("G37").Value = "13"
$ws.Range("G38").Value = "13"
$ws.Range("G39").Value = "13"
$ws.Range("D40").Value = "0.03987"
$ws.Range("G40").Value = "13"
$ws.Range("D41").Value = "0.006512"
$ws.Range("G41").Value = "13"
$ws.Range("D42").Value = "0.1071"
$ws.Range("G42").Value = "13"
$ws.Range("D43").Value = "0.003201"
$ws.Range("G43").Value = "13"
$ws.Range("D44").Value = "0.005809"
$ws.Range("G44").Value = "13"
$ws.Range("D45").Value = "0.00005277"
$ws.Range("G45").Value = "13"
$ws.Range("G46").Value = "13"
$ws.Range("D47").Value = "0.9903"
$ws.Range("G47").Value = "13"
$ws.Range("D48").Value = "0.01449"
$ws.Range("G48").Value = "13"
$ws.Range("G49").Value = "13"
$ws.Range("G50").Value = "13"
$ws.Range("G51").Value = "13"

# --- Plain text (Volume/coin label) column (E) updates ---
$ws.Range("E17").Value = "16OneONE"
$ws.Range("E27").Value = "26UpBotsUBXT"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
